$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 436.66666
$ws.Range("J97").Value = 436.66666
$ws.Range("L97").Value = 1309.99998
$ws.Range("N97").Value = -2301.99998
$ws.Range("H101").Value = 822.75
$ws.Range("I101").Value = 897.4286
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 2692.2858
$ws.Range("L101").Value = 900
$ws.Range("M101").Value = -1070.2858
$ws.Range("N101").Value = -4144
$ws.Range("H118").Value = 2008.3334
$ws.Range("I118").Value = 386.66666
$ws.Range("J118").Value = 3630
$ws.Range("K118").Value = 1159.99998
$ws.Range("L118").Value = 10890
$ws.Range("M118").Value = 497.0000199999999
$ws.Range("N118").Value = -14204
$ws.Range("H127").Value = 359972.3
$ws.Range("I127").Value = 396.2
$ws.Range("J127").Value = 431887.53
$ws.Range("K127").Value = 1188.6
$ws.Range("L127").Value = 1295662.59
$ws.Range("M127").Value = 3771.4
$ws.Range("N127").Value = -1305582.59
$ws.Range("H129").Value = 1300.5834
$ws.Range("I129").Value = 469.9
$ws.Range("J129").Value = 1412.8379
$ws.Range("K129").Value = 1409.7
$ws.Range("L129").Value = 4238.5137
$ws.Range("M129").Value = 3590.3
$ws.Range("N129").Value = -14238.5137
$ws.Range("H132").Value = 6253380
$ws.Range("I132").Value = 2767.6128
$ws.Range("J132").Value = 27783266
$ws.Range("K132").Value = 8302.8384
$ws.Range("L132").Value = 83349798
$ws.Range("M132").Value = -5772.838400000001
$ws.Range("N132").Value = -83354858
$ws.Range("H138").Value = 9261863
$ws.Range("I138").Value = 15875118
$ws.Range("J138").Value = 3306.9333
$ws.Range("K138").Value = 47625354
$ws.Range("L138").Value = 9920.7999
$ws.Range("M138").Value = -47620214
$ws.Range("N138").Value = -20200.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10361.506
$ws.Range("I32").Value = 9577.411
$ws.Range("J32").Value = 13938.9375
$ws.Range("K32").Value = 9577.411
$ws.Range("L32").Value = 13938.9375
$ws.Range("M32").Value = -9290.411
$ws.Range("N32").Value = -14512.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3563.0232
$ws.Range("I134").Value = 2463.8823
$ws.Range("J134").Value = 7715.3335
$ws.Range("K134").Value = 7391.646900000001
$ws.Range("L134").Value = 23146.0005
$ws.Range("M134").Value = -4856.646900000001
$ws.Range("N134").Value = -28216.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7756181.5
$ws.Range("I31").Value = 4470.4707
$ws.Range("J31").Value = 37040424
$ws.Range("K31").Value = 4470.4707
$ws.Range("L31").Value = 37040424
$ws.Range("M31").Value = -4175.4707
$ws.Range("N31").Value = -37041014
$ws.Range("H34").Value = 7756181.5
$ws.Range("I34").Value = 4470.4707
$ws.Range("J34").Value = 37040424
$ws.Range("K34").Value = 4470.4707
$ws.Range("L34").Value = 37040424
$ws.Range("M34").Value = -4268.4707
$ws.Range("N34").Value = -37040828

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 5517.5454
$ws.Range("I125").Value = 2030
$ws.Range("J125").Value = 6292.5557
$ws.Range("K125").Value = 6090
$ws.Range("L125").Value = 18877.6671
$ws.Range("M125").Value = -1170
$ws.Range("N125").Value = -28717.6671
$ws.Range("H131").Value = 1492.0405
$ws.Range("I131").Value = 2289.45
$ws.Range("J131").Value = 1196.7037
$ws.Range("K131").Value = 6868.349999999999
$ws.Range("L131").Value = 3590.1111
$ws.Range("M131").Value = -1828.349999999999
$ws.Range("N131").Value = -13670.1111
$ws.Range("H133").Value = 40003964
$ws.Range("I133").Value = 76924970
$ws.Range("J133").Value = 6215.8335
$ws.Range("K133").Value = 230774910
$ws.Range("L133").Value = 18647.5005
$ws.Range("M133").Value = -230769850
$ws.Range("N133").Value = -28767.5005
$ws.Range("H134").Value = 4860.375
$ws.Range("I134").Value = 2013.5454
$ws.Range("J134").Value = 7269.231
$ws.Range("K134").Value = 6040.6362
$ws.Range("L134").Value = 21807.693
$ws.Range("M134").Value = -970.6361999999999
$ws.Range("N134").Value = -31947.693
$ws.Range("H136").Value = 3394.1482
$ws.Range("I136").Value = 991.7273
$ws.Range("J136").Value = 5045.8125
$ws.Range("K136").Value = 2975.1819
$ws.Range("L136").Value = 15137.4375
$ws.Range("M136").Value = 2124.8181
$ws.Range("N136").Value = -25337.4375
$ws.Range("H137").Value = 4739.3213
$ws.Range("I137").Value = 2641.75
$ws.Range("J137").Value = 7536.0835
$ws.Range("K137").Value = 7925.25
$ws.Range("L137").Value = 22608.2505
$ws.Range("M137").Value = -2825.25
$ws.Range("N137").Value = -32808.25049999999
$ws.Range("H139").Value = 2348.2917
$ws.Range("I139").Value = 1326.619
$ws.Range("K139").Value = 3979.857
$ws.Range("M139").Value = 1160.143
$ws.Range("H140").Value = 3967.4583
$ws.Range("I140").Value = 2214.6
$ws.Range("J140").Value = 6888.8887
$ws.Range("K140").Value = 6643.799999999999
$ws.Range("L140").Value = 20666.6661
$ws.Range("M140").Value = -1463.799999999999
$ws.Range("N140").Value = -31026.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4418.1113
$ws.Range("I126").Value = 3193.2727
$ws.Range("J126").Value = 6342.857
$ws.Range("K126").Value = 9579.8181
$ws.Range("L126").Value = 19028.571
$ws.Range("M126").Value = -7109.8181
$ws.Range("N126").Value = -23968.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5868.815
$ws.Range("I7").Value = 8065
$ws.Range("J7").Value = 4770.722
$ws.Range("K7").Value = 8065
$ws.Range("L7").Value = 4770.722
$ws.Range("M7").Value = -7953
$ws.Range("N7").Value = -4994.722
$ws.Range("H126").Value = 5868.815
$ws.Range("I126").Value = 8065
$ws.Range("J126").Value = 4770.722
$ws.Range("K126").Value = 24195
$ws.Range("L126").Value = 14312.166
$ws.Range("M126").Value = -21725
$ws.Range("N126").Value = -19252.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6624.5
$ws.Range("J62").Value = 6799.6
$ws.Range("L62").Value = 6799.6
$ws.Range("N62").Value = -8047.6
$ws.Range("H65").Value = 6624.5
$ws.Range("J65").Value = 6799.6
$ws.Range("L65").Value = 33998
$ws.Range("N65").Value = -40238
$ws.Range("H132").Value = 1798.8889
$ws.Range("I132").Value = 1749.5555
$ws.Range("J132").Value = 1897.5555
$ws.Range("K132").Value = 5248.666499999999
$ws.Range("L132").Value = 5692.666499999999
$ws.Range("M132").Value = -2718.666499999999
$ws.Range("N132").Value = -10752.6665
